$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the date number format (same format code previously used for E3:E5)
# and set the date values (44784 = 2022-08-11) on those cells.
$ws.Range("E3:E5").NumberFormat = "[$-409]mmm\-yy;@"
$ws.Range("E3").Value = 44784
$ws.Range("E4").Value = 44784
$ws.Range("E5").Value = 44784

# New trailing cell with a new shared string "SS".
$ws.Range("G6").Value = "SS"

# Move/update the active selection to the newly written cell, matching the
# new used range (dimension grows to A1:G6).
$ws.Range("G6").Select() | Out-Null
